$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record was inserted at row 476, shifting the existing
# rows 476-506 down by one (to 477-507).
$ws.Rows.Item(476).Insert()

# Populate the newly inserted row 476 with the new record's data.
$ws.Cells.Item(476, 1).Value2 = 6
$ws.Cells.Item(476, 2).Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(476, 3).Value2 = "Metropolitana"
$ws.Cells.Item(476, 4).Value2 = 45267
$ws.Cells.Item(476, 5).Value2 = 13
$ws.Cells.Item(476, 6).Value2 = 100112026
$ws.Cells.Item(476, 7).Value2 = "Haba"
$ws.Cells.Item(476, 8).Value2 = "Sin especificar"
$ws.Cells.Item(476, 9).Value2 = "Primera"
$ws.Cells.Item(476, 10).Value2 = 500
$ws.Cells.Item(476, 11).Value2 = 11000
$ws.Cells.Item(476, 12).Value2 = 13000
$ws.Cells.Item(476, 13).Value2 = 11920
$ws.Cells.Item(476, 14).Value2 = "`$/saco 25 kilos"
$ws.Cells.Item(476, 15).Value2 = "Región del Maule"
$ws.Cells.Item(476, 16).Value2 = 477
$ws.Cells.Item(476, 17).Value2 = 25
$ws.Cells.Item(476, 18).Value2 = "Hortaliza"
